# Q3 Update - 2025
# Applies the UNHCR "fromCSV" refresh:
#   1. short-url column (B2:B629) refreshed from "68SjwG" to "2mK7Ae"
#   2. Row 615 country-of-origin corrected from Niger (NGR/NER, coo_id 139)
#      to Nigeria (NIG/NGA, coo_id 141)
#   3. A handful of numeric statistics cells (N/O/P/S/V) refreshed with the
#      latest quarterly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the short-url value for every data row (row 2 .. row 629).
#    The sheet stores this as a shared string, so updating every cell to
#    the same literal collapses back down to a single shared entry on
#    save, exactly mirroring the source diff.
# ---------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 629 }
$ws.Range("B2:B$lastRow").Value = "2mK7Ae"

# ---------------------------------------------------------------------
# 2) Numeric-looking statistic cells must stay stored as text, matching
#    how this CSV-imported sheet already represents every number. Excel
#    auto-coerces a digit-only string into a real number on assignment,
#    so mark the target cells as Text ("@") first, then write the new
#    figures.
# ---------------------------------------------------------------------
$numericTextCells = @(
    "N602","O602",
    "N605",
    "O607",
    "N608",
    "N613",
    "F615",
    "N616","O616",
    "N617","O617",
    "N618",
    "N620","O620",
    "N621","O621","P621",
    "S622",
    "N623","O623",
    "N624","O624",
    "N627",
    "V628",
    "N629"
)
$numericTextRange = $ws.Range($numericTextCells[0])
foreach ($addr in $numericTextCells) {
    $numericTextRange = $excel.Union($numericTextRange, $ws.Range($addr))
}
$numericTextRange.NumberFormat = "@"

$ws.Range("N602").Value = "234"
$ws.Range("O602").Value = "81"
$ws.Range("N605").Value = "28"
$ws.Range("O607").Value = "6"
$ws.Range("N608").Value = "40"
$ws.Range("N613").Value = "231"
$ws.Range("F615").Value = "141"
$ws.Range("N616").Value = "648"
$ws.Range("O616").Value = "156"
$ws.Range("N617").Value = "22393"
$ws.Range("O617").Value = "1822"
$ws.Range("N618").Value = "5"
$ws.Range("N620").Value = "49122"
$ws.Range("O620").Value = "1710"
$ws.Range("N621").Value = "975020"
$ws.Range("O621").Value = "59"
$ws.Range("P621").Value = "12874"
$ws.Range("S622").Value = "10284"
$ws.Range("N623").Value = "65367"
$ws.Range("O623").Value = "168"
$ws.Range("N624").Value = "42"
$ws.Range("O624").Value = "5"
$ws.Range("N627").Value = "95"
$ws.Range("V628").Value = "3200000"
$ws.Range("N629").Value = "40"

# ---------------------------------------------------------------------
# 3) Row 615 country-of-origin correction: Niger (139/NGR/NER) was
#    actually Nigeria (141/NIG/NGA).
# ---------------------------------------------------------------------
$ws.Range("G615").Value = "Nigeria"
$ws.Range("H615").Value = "NIG"
$ws.Range("I615").Value = "NGA"
